$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33..111 down to 34..112
# (Excel's default Insert shifts cells down and copies formatting from the row above,
# which preserves the date-format style "s=2" already on column D).
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly price record.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 45012
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = "Tropicales y subtropicales"
$ws.Range("I33").Value = 100108002
$ws.Range("J33").Value = "Mango"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 7000
$ws.Range("O33").Value = 7000
$ws.Range("P33").Value = 7000
$ws.Range("Q33").Value = '$/bandeja 4 kilos'
$ws.Range("R33").Value = "Perú"
$ws.Range("S33").Value = 1750
$ws.Range("T33").Value = 4
